$wb = $excel.ActiveWorkbook

# ===== numeric_variable =====
$ws = $wb.Worksheets.Item("numeric_variable")
$ws.Range("N2").Value = 0
$ws.Range("Q2").Value = 1
$ws.Range("T2").Value = 2
$ws.Range("N4").Value = 0.7041217312852566
$ws.Range("O4").Value = -1
$ws.Range("P4").Value = 0.4082434625705131
$ws.Range("Q4").Value = 0.5497199333020351
$ws.Range("S4").Value = 0.09943986660407012
$ws.Range("T4").Value = 0.7741065466434253
$ws.Range("U4").Value = 1
$ws.Range("V4").Value = 0.5482130932868505
$ws.Range("B15").Value = 9.869827596845477
$ws.Range("B16").Value = 4.951099843665786
$ws.Range("B17").Value = -9.458561131516781
$ws.Range("B18").Value = 6.552760643691096
$ws.Range("B19").Value = 9.988824899049497
$ws.Range("B20").Value = 13.28321426209553
$ws.Range("B21").Value = 23.70176620263562

# ===== object_variable =====
$ws = $wb.Worksheets.Item("object_variable")
# Rotate D13/F13/H13 text labels via copy/paste (avoids Excel auto-converting "0.0%" etc. to a number)
$ws.Range("H13").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4104) | Out-Null
$ws.Range("F13").Copy() | Out-Null
$ws.Range("H13").PasteSpecial(-4104) | Out-Null
$ws.Range("D13").Copy() | Out-Null
$ws.Range("F13").PasteSpecial(-4104) | Out-Null
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4104) | Out-Null
$ws.Range("Z1").Clear() | Out-Null
$excel.CutCopyMode = $false

$ws.Range("N2").Value = 0
$ws.Range("Q2").Value = 1
$ws.Range("T2").Value = 2
$ws.Range("N4").Value = 0.9283940455097172
$ws.Range("P4").Value = 0.8567880910194343
$ws.Range("Q4").Value = 0.8249987858380147
$ws.Range("S4").Value = 0.6499975716760293
$ws.Range("T4").Value = 0.9276348310221921
$ws.Range("V4").Value = 0.8552696620443843
$ws.Range("C13").Value = 0
$ws.Range("E13").Value = 1
$ws.Range("G13").Value = 2
$ws.Range("A14").Value = 5
$ws.Range("B14").Value = 117
$ws.Range("C14").Value = 15
$ws.Range("D14").Value = 12.82051282051282
$ws.Range("E14").Value = 101
$ws.Range("F14").Value = 86.32478632478633
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 0.8547008547008547
$ws.Range("A15").Value = 6
$ws.Range("B15").Value = 112
$ws.Range("C15").Value = 9
$ws.Range("D15").Value = 8.035714285714286
$ws.Range("E15").Value = 97
$ws.Range("F15").Value = 86.60714285714286
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = 5.357142857142857
$ws.Range("A16").Value = 2
$ws.Range("B16").Value = 103
$ws.Range("D16").Value = 76.69902912621359
$ws.Range("E16").Value = 24
$ws.Range("F16").Value = 23.30097087378641
$ws.Range("A17").Value = 8
$ws.Range("B17").Value = 100
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 70
$ws.Range("F17").Value = 70
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = 30
$ws.Range("A18").Value = 9
$ws.Range("B18").Value = 94
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 51
$ws.Range("F18").Value = 54.25531914893617
$ws.Range("G18").Value = 43
$ws.Range("H18").Value = 45.74468085106383
$ws.Range("A19").Value = 3
$ws.Range("B19").Value = 94
$ws.Range("C19").Value = 50
$ws.Range("D19").Value = 53.19148936170212
$ws.Range("E19").Value = 44
$ws.Range("F19").Value = 46.80851063829788
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = 91
$ws.Range("C20").Value = 30
$ws.Range("D20").Value = 32.96703296703296
$ws.Range("E20").Value = 61
$ws.Range("F20").Value = 67.03296703296704
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("A21").Value = 7
$ws.Range("B21").Value = 91
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 80
$ws.Range("F21").Value = 87.91208791208791
$ws.Range("G21").Value = 11
$ws.Range("H21").Value = 12.08791208791209
$ws.Range("A22").Value = 1
$ws.Range("B22").Value = 90
$ws.Range("C22").Value = 76
$ws.Range("D22").Value = 84.44444444444444
$ws.Range("E22").Value = 14
$ws.Range("F22").Value = 15.55555555555556
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("B23").Value = 57
$ws.Range("C23").Value = 54
$ws.Range("D23").Value = 94.73684210526316
$ws.Range("E23").Value = 3
$ws.Range("F23").Value = 5.263157894736842
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("B24").Value = 51
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = 17.64705882352941
$ws.Range("G24").Value = 42
$ws.Range("H24").Value = 82.35294117647059

# ===== column with soo000oo000oo000oo0 =====
$ws = $wb.Worksheets.Item("column with soo000oo000oo000oo0")
$ws.Range("N2").Value = 0
$ws.Range("Q2").Value = 1
$ws.Range("T2").Value = 2
$ws.Range("N4").Value = 0.5009370741893029
$ws.Range("P4").Value = 0.001874148378605778
$ws.Range("Q4").Value = 0.5126353790613718
$ws.Range("R4").Value = -1
$ws.Range("S4").Value = 0.0252707581227436
$ws.Range("T4").Value = 0.5288220551378446
$ws.Range("V4").Value = 0.05764411027568928
$ws.Range("B15").Value = 9.713528973077924
$ws.Range("B16").Value = 5.023536697732114
$ws.Range("B17").Value = -5.860648175735951
$ws.Range("B18").Value = 6.149397526742504
$ws.Range("B19").Value = 9.718042688703704
$ws.Range("B20").Value = 13.0453282926783
$ws.Range("B21").Value = 24.82810229050938
